$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "Normalized Weight" column (B) next to the existing
# "Peronal Weight" column (A) in each of the three personal-score blocks
# (rows 2, 12, 22), and the matching column (J) next to the second
# "Peronal Weight" column (I) in the first two blocks (rows 2, 12).
# The new cells reuse the same rotated/centered style as their neighbor.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Normalized Weight"
$ws.Range("J2").Value = "Normalized Weight"
$ws.Range("B12").Value = "Normalized Weight"
$ws.Range("J12").Value = "Normalized Weight"
$ws.Range("B22").Value = "Normalized Weight"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null

$ws.Range("A12").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Copy() | Out-Null
$ws.Range("J12").PasteSpecial(-4122) | Out-Null

$ws.Range("A22").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Third judge's raw scores (rows 13-20) were filled in (they previously had
# gaps), which feeds the "Normalized Weight" calculations in column J and
# the final weighted totals lower in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("I13").Value = 2
$ws.Range("L13").Value = 10
$ws.Range("M13").Value = 8
$ws.Range("N13").Value = 10
$ws.Range("O13").Value = 8

$ws.Range("I14").Value = 8
$ws.Range("L14").Value = 2
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 7

$ws.Range("I15").Value = 8
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 10
$ws.Range("N15").Value = 4
$ws.Range("O15").Value = 7

$ws.Range("I16").Value = 7
$ws.Range("L16").Value = 6
$ws.Range("M16").Value = 9
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 4

$ws.Range("I17").Value = 5
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 10
$ws.Range("N17").Value = 4
$ws.Range("O17").Value = 6

$ws.Range("I18").Value = 9
$ws.Range("L18").Value = 5
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 6

$ws.Range("I19").Value = 10
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 8
$ws.Range("N19").Value = 3
$ws.Range("O19").Value = 4

$ws.Range("I20").Value = 7
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = 7
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 3

# Leave the selection where the author finished working, like the saved
# workbook records.
$ws.Range("B22").Select() | Out-Null

$wb.Save()
